$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row in column A (currently A1:A13 are populated -> new row 14)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Add the new vendor entry
$ws.Cells.Item($newRow, 1).Value = "zuluCrypt"

# Match formatting used by the previous last row (A13)
$ws.Cells.Item($lastRow, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Keep the active selection on the newly added cell, like in the source workbook
$ws.Range("A" + $newRow).Select()
